$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 06:50"

# Row 4: Estados Unidos (USA) - update totals
$ws.Range("B4").Value = 215300
$ws.Range("C4").Value = 297
$ws.Range("E4").Value = 201312

# Row 37: Pakistan - update totals
$ws.Range("B37").Value = 2238
$ws.Range("C37").Value = 120
$ws.Range("E37").Value = 2113
$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 31

# Rows 115-116: swap Bolivia and Kirguistan (country name + stats)
# Row 115 becomes Kirguistan with updated stats
$ws.Range("A115").Value = "Kirguistan"
$ws.Range("B115").Value = 116
$ws.Range("C115").Value = 5
$ws.Range("D115").Value = 5
$ws.Range("E115").Value = 111
$ws.Range("F115").Value = 5
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0

# Row 116 becomes Bolivia, carrying forward its previous (unchanged) stats
$ws.Range("A116").Value = "Bolivia"
$ws.Range("B116").Value = 115
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 1
$ws.Range("E116").Value = 107
$ws.Range("F116").Value = 3
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 7
